# Commit: "regen save_data to use K instead of Strike#, regen std/mean,
# calc and write s_vals"
#
# The underlying save-data pipeline recomputed the "K" column (header in
# G1) for every saved trade row and rewrote the resulting integer values
# back into the sheet. Reproduce that by writing the newly calculated
# values straight into column G for each data row (rows 2-59).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new K value
$newK = [ordered]@{
    2  = 1
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 4
    8  = 0
    9  = 0
    10 = 1
    11 = 0
    12 = 2
    13 = 0
    14 = 1
    15 = 2
    16 = 4
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 2
    25 = 2
    26 = 1
    27 = 4
    28 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 0
    36 = 1
    37 = 2
    38 = 1
    39 = 2
    40 = 1
    41 = 0
    42 = 4
    43 = 0
    44 = 0
    45 = 1
    46 = 0
    47 = 2
    48 = 2
    49 = 1
    50 = 3
    51 = 1
    52 = 1
    53 = 0
    54 = 2
    55 = 0
    56 = 1
    57 = 1
    58 = 2
    59 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
